$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("phylip-programs")
$ws1.Activate() | Out-Null

$ws1.Range("C3").Copy($ws1.Range("C22"))
$ws1.Range("C3").Copy($ws1.Range("C25"))
$ws1.Range("C3").Copy($ws1.Range("C26"))

$ws1.Range("A22").Copy($ws1.Range("B22"))
$ws1.Range("A22").Copy($ws1.Range("B25"))
$ws1.Range("A22").Copy($ws1.Range("B26"))

$ws1.Range("B25").Value = "Rpars"
$ws1.Range("B22").Value = "Rmix"
$ws1.Range("B26").Value = "Rpenny"

$ws1.Range("C22").Value = 41638
$ws1.Range("C25").Value = 41634
$ws1.Range("C26").Value = 41638

$ws1.Range("B2").Select() | Out-Null
